$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @{ B = 0.2917716402565462;  C = 0.002571899574220771; D = 0.7527432677738641; E = 0.4942365360607697; G = 1.541323343665401 }
    3 = @{ B = 3.286832544864788;   C = 1.655778082260271;    D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    4 = @{ B = 3.286832544864788;   C = 1.655778082260271;    D = 0.1494219747398047; E = 0.4942365360607697; G = 5.586269137925634 }
    5 = @{ B = 3.286832544864788;   C = 0.306821227259698;    D = 0.1494219747398047; E = 0.4942365360607697; G = 4.23731228292506 }
    6 = @{ B = 3.286832544864788;   C = 1.655778082260271;    D = 0.7527432677738641; E = 0.4942365360607697; G = 6.189590430959694 }
    7 = @{ B = 3.286832544864788;   C = 1.655778082260271;    D = 22.3905356188092;   E = 0.4942365360607697; G = 27.82738278199502 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
